# The workbook originally has 5 columns (A:E) and 3 data rows (1:3):
#   Row, 1-c__Fusobacteriia, max, prediction, rejection-f
#   even_MAG-GUT88654.fa, 1, 1, c__Fusobacteriia, c__Fusobacteriia
#   even_MAG-GUT88709.fa, 1, 1, c__Fusobacteriia, c__Fusobacteriia
#
# The target output keeps a single child/prediction row and drops the
# "max" column entirely, shrinking the sheet to A1:D2:
#   Row, 1-c__Fusobacteriia, prediction, rejection-f
#   even_MAG-GUT88654.fa, 10990.46771063232, c__Fusobacteriia, c__Fusobacteriia

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column C ("max"); columns D and E shift left into C and D.
$ws.Columns.Item(3).Delete()

# Remove the now second data row (former row 3, "even_MAG-GUT88709.fa").
$ws.Rows.Item(3).Delete()

# Update the remaining data row's value cell with the new numeric value.
$ws.Range("B2").Value = 10990.46771063232
